# Updated symbol list on Mon Jan 30 19:27:46 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns D (Price) and E (Volume(1h)) hold numeric-looking values that are
# stored as plain text in the workbook. Prefix with a leading apostrophe so
# Excel keeps them as text instead of auto-converting to numbers/percentages.

$ws.Range("D2").Value = '''306.65'
$ws.Range("E2").Value = '''-3.76%'
$ws.Range("D3").Value = '''37.33'
$ws.Range("E3").Value = '''-6.17%'
$ws.Range("D4").Value = '''5.101'
$ws.Range("E4").Value = '''-0.92%'
$ws.Range("D5").Value = '''0.07748'
$ws.Range("E5").Value = '''-5.83%'
$ws.Range("B6").Value = 'GateToken'
$ws.Range("C6").Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range("D6").Value = '''4.393'
$ws.Range("E6").Value = '''1.63%'
$ws.Range("B7").Value = 'FTXToken'
$ws.Range("C7").Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range("D7").Value = '''1.916'
$ws.Range("E7").Value = '''-8.26%'
$ws.Range("D8").Value = '''8.214'
$ws.Range("E8").Value = '''-1.22%'
$ws.Range("D9").Value = '''3.085'
$ws.Range("E9").Value = '''-8.33%'
$ws.Range("D10").Value = '''0.9251'
$ws.Range("E10").Value = '''-1.67%'
$ws.Range("D11").Value = '''0.1277'
$ws.Range("E11").Value = '''-6.31%'
$ws.Range("D12").Value = '''0.1899'
$ws.Range("E12").Value = '''-4.41%'
$ws.Range("D13").Value = '''0.08806'
$ws.Range("E13").Value = '''-3.56%'
$ws.Range("D14").Value = '''0.03434'
$ws.Range("E14").Value = '''-1.68%'
$ws.Range("D15").Value = '''0.09719'
$ws.Range("E15").Value = '''-1.18%'
$ws.Range("D16").Value = '''0.001368'
$ws.Range("E16").Value = '''-2.26%'
$ws.Range("D17").Value = '''0.005937'
$ws.Range("E17").Value = '''-6.00%'
$ws.Range("D18").Value = '''3.597'
$ws.Range("E18").Value = '''-2.66%'
$ws.Range("D19").Value = '''0.3386'
$ws.Range("E19").Value = '''-2.65%'
$ws.Range("E20").Value = '''-1.67%'
$ws.Range("D21").Value = '''5.042'
$ws.Range("E21").Value = '''0.74%'
$ws.Range("D22").Value = '''0.2500'
$ws.Range("E22").Value = '''2.28%'
$ws.Range("D23").Value = '''0.02120'
$ws.Range("E23").Value = '''5,213.63%'
$ws.Range("D24").Value = '''0.04354'
$ws.Range("E24").Value = '''0.34%'
$ws.Range("D25").Value = '''0.001225'
$ws.Range("E25").Value = '''0.03%'
$ws.Range("D26").Value = '''0.004506'
$ws.Range("E26").Value = '''-6.64%'
$ws.Range("D27").Value = '''0.0001361'
$ws.Range("E27").Value = '''5.06%'
$ws.Range("D39").Value = '''0.02201'
$ws.Range("E39").Value = '''-0.92%'
$ws.Range("D40").Value = '''0.04929'
$ws.Range("E40").Value = '''-5.46%'
$ws.Range("D41").Value = '''0.007731'
$ws.Range("E41").Value = '''0.69%'
$ws.Range("D42").Value = '''0.009816'
$ws.Range("E42").Value = '''1.76%'
$ws.Range("E43").Value = '''-5.00%'
$ws.Range("D44").Value = '''0.002009'
$ws.Range("E44").Value = '''-5.35%'
$ws.Range("D45").Value = '''0.008392'
$ws.Range("E45").Value = '''-6.08%'
$ws.Range("D46").Value = '''0.00006861'
$ws.Range("E46").Value = '''3.16%'
$ws.Range("D47").Value = '''0.00000000754'
$ws.Range("E47").Value = '''0.82%'
$ws.Range("D48").Value = '''0.003020'
$ws.Range("E48").Value = '''4.98%'
$ws.Range("D49").Value = '''0.001306'
$ws.Range("E49").Value = '''-22.49%'
$ws.Range("D50").Value = '''0.00002110'
$ws.Range("E50").Value = '''0.82%'
$ws.Range("D51").Value = '''0.0002009'
$ws.Range("E51").Value = '''0.82%'
